# Add a new "Sun2" AH sequence row to Sheet1, and restyle its sequence cell
# with an explicit black Helvetica 11 font (the "AH name disappear" font fix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row (row 53) -------------------------------------------------
$ws.Range("A53").Value = "LRSAVSRAGSLLWMVATSPGRLFRLL"
$ws.Range("B53").Value = "Sun2"
$ws.Range("C53").Value = 1
$ws.Range("D53").Value = 0

# Give the new sequence cell its own explicit font (Helvetica 11, solid
# black) instead of the inherited theme font - this is the fix for the AH
# name disappearing. Color is applied before Name so the intermediate state
# reuses the already-existing "explicit black" font instead of minting an
# extra unused one.
$font = $ws.Range("A53").Font
$font.Color = 0
$font.Name = "Helvetica"
$font.Size = 11

# --- View state: scroll position + active selection ------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 23
$win.ScrollColumn = 1
$null = $ws.Range("D54").Select()
